# Daily attendance processing - 2025-11-02 21:18:56
#
# Column G ("Recorded By") holds a comma-separated list of the users who
# recorded/edited a session's attendance. This pass re-normalises the
# display order of that list: the whole list is reversed, except that a
# literal lowercase "system" entry (the automated/bot actor) is pinned in
# its original position while everything else around it reverses.
#
# NOTE: this interpreter's PowerShell functions do not get their own
# variable scope - a loop counter reused inside a helper function clobbers
# the same-named counter in the caller. Every loop/temp variable below is
# therefore given a unique name to avoid cross-talk between the helper
# functions and their callers.

function Test-ExactEqual($exEqA, $exEqB) {
    if ($exEqA.Length -ne $exEqB.Length) { return $false }
    $exEqCharsA = $exEqA.ToCharArray()
    $exEqCharsB = $exEqB.ToCharArray()
    for ($exEqIdx = 0; $exEqIdx -lt $exEqCharsA.Length; $exEqIdx++) {
        if ([int]$exEqCharsA[$exEqIdx] -ne [int]$exEqCharsB[$exEqIdx]) { return $false }
    }
    return $true
}

function Transform-RecordedBy($trVal) {
    $trParts = $trVal.Split(",")
    for ($trTrimIdx = 0; $trTrimIdx -lt $trParts.Count; $trTrimIdx++) {
        $trParts[$trTrimIdx] = $trParts[$trTrimIdx].Trim()
    }
    $trN = $trParts.Count

    # indices of entries that are NOT the pinned literal "system" token
    $trNonPinIdx = @()
    foreach ($trScanIdx in 0..($trN - 1)) {
        $trIsPin = Test-ExactEqual $trParts[$trScanIdx] "system"
        if (-not $trIsPin) {
            $trNonPinIdx += $trScanIdx
        }
    }

    $trResult = New-Object 'object[]' $trN
    for ($trCopyIdx = 0; $trCopyIdx -lt $trN; $trCopyIdx++) { $trResult[$trCopyIdx] = $trParts[$trCopyIdx] }

    $trM = $trNonPinIdx.Count
    for ($trSwapIdx = 0; $trSwapIdx -lt $trM; $trSwapIdx++) {
        $trSrcIdx = $trNonPinIdx[$trM - 1 - $trSwapIdx]
        $trDstIdx = $trNonPinIdx[$trSwapIdx]
        $trResult[$trDstIdx] = $trParts[$trSrcIdx]
    }

    return ($trResult -join ", ")
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Column G is the 7th column ("Recorded By")
$recordedByCol = 7

for ($rowIdx = $firstRow; $rowIdx -le $lastRow; $rowIdx++) {
    $cell = $ws.Cells.Item($rowIdx, $recordedByCol)
    $origText = $cell.Text
    if ($origText -ne $null -and $origText.Length -gt 0) {
        $newText = Transform-RecordedBy $origText
        if ($newText -ne $origText) {
            $cell.Value = $newText
        }
    }
}
